$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" note cell (A1) with the refreshed conversion text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.87 = 10786.22 pesos`n✅ 10786.22 pesos = 2.86 = 954.91 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet numeric cells ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("O10").Value = 3756
$wsTasas.Range("N12").Value = 3777.24
$wsTasas.Range("O12").Value = 334.4
